# casos_pais.xlsx: add the EEUU (USA) column and fill in the row that was
# missing its "regional" link data (row 41 / the 26th day), per the commit
# "Me faltó agregar el link al regional del 26".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column header -----------------------------------------------------
$ws.Range("H1").Value = "EEUU"

# --- Row 41 was missing its Chile/Brasil/España/Italia/Corea.del.Sur values ---
$ws.Range("C41").Value = 1306
$ws.Range("D41").Value = 2554
$ws.Range("E41").Value = 49515
$ws.Range("F41").Value = 74386
$ws.Range("G41").Value = 9137

# --- New EEUU data for rows 2..41 ------------------------------------------
$eeuu = @(15,15,15,15,15,15,35,35,35,53,57,60,60,63,68,75,100,124,158,221,319,435,541,704,994,1301,1630,2183,2770,3613,4596,6344,9197,13779,19367,24192,33592,43781,54856,68211)
for ($i = 0; $i -lt $eeuu.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 8).Value = $eeuu[$i]
}

# --- Reflect that the user scrolled down to row 26 while fixing it ---------
$excel.ActiveWindow.FreezePanes = $false
[void]$ws.Range("B26").Select()
$excel.ActiveWindow.FreezePanes = $true
[void]$ws.Range("H2").Select()
